$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$orig_D2 = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.239.10'
$ws.Range('D2').Style = $orig_D2
$ws.Range('E2').Value = '  -0.65%  '
$orig_D3 = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.839.25'
$ws.Range('D3').Style = $orig_D3
$ws.Range('E3').Value = '  -1.42%  '
$orig_D4 = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9999'
$ws.Range('D4').Style = $orig_D4
$ws.Range('E4').Value = '  -0.06%  '
$orig_D5 = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.17'
$ws.Range('D5').Style = $orig_D5
$ws.Range('E5').Value = '  -1.07%  '
$ws.Range('E6').Value = '  +0.00%  '
$orig_D7 = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4656'
$ws.Range('D7').Style = $orig_D7
$ws.Range('E7').Value = '  -2.61%  '
$ws.Range('E8').Value = '  -2.37%  '
$ws.Range('E9').Value = '  -4.07%  '
$orig_D10 = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.838.85'
$ws.Range('D10').Style = $orig_D10
$ws.Range('E10').Value = '  -1.22%  '
$orig_D11 = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07428'
$ws.Range('D11').Style = $orig_D11
$ws.Range('E11').Value = '  -0.24%  '
$orig_D12 = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.21'
$ws.Range('D12').Style = $orig_D12
$ws.Range('E12').Value = '  +0.00%  '
$orig_D13 = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.928'
$ws.Range('D13').Style = $orig_D13
$ws.Range('E13').Value = '  -2.59%  '
$orig_D14 = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '83.49'
$ws.Range('D14').Style = $orig_D14
$ws.Range('E14').Value = '  -3.92%  '
$orig_D15 = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6190'
$ws.Range('D15').Style = $orig_D15
$ws.Range('E15').Value = '  -3.15%  '
$orig_D16 = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.177.84'
$ws.Range('D16').Style = $orig_D16
$ws.Range('E16').Value = '  -0.82%  '
$orig_D18 = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '227.19'
$ws.Range('D18').Style = $orig_D18
$ws.Range('E18').Value = '  -2.86%  '
$orig_D19 = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007281'
$ws.Range('D19').Style = $orig_D19
$ws.Range('E19').Value = '  -2.41%  '
$ws.Range('E20').Value = '  -4.99%  '
$orig_D21 = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').Style = $orig_D21
$ws.Range('E21').Value = '  +0.03%  '
$orig_D22 = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.887'
$ws.Range('D22').Style = $orig_D22
$ws.Range('E22').Value = '  -4.56%  '
$orig_D23 = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.852'
$ws.Range('D23').Style = $orig_D23
$ws.Range('E23').Value = '  -3.89%  '
$orig_D24 = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.176'
$ws.Range('D24').Style = $orig_D24
$ws.Range('E24').Value = '  -1.34%  '
$orig_D25 = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '164.19'
$ws.Range('D25').Style = $orig_D25
$ws.Range('E25').Value = '  -2.72%  '
$ws.Range('E26').Value = '  -2.16%  '
$orig_D27 = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.866'
$ws.Range('D27').Style = $orig_D27
$ws.Range('E27').Value = '  -1.54%  '
$orig_D28 = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1032'
$ws.Range('D28').Style = $orig_D28
$ws.Range('E28').Value = '  -2.17%  '
$orig_D29 = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.372'
$ws.Range('D29').Style = $orig_D29
$ws.Range('E29').Value = '  -0.51%  '
$orig_D30 = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.076'
$ws.Range('D30').Style = $orig_D30
$ws.Range('E30').Value = '  -4.48%  '
$ws.Range('E31').Value = '  -4.27%  '
$orig_D32 = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.04832'
$ws.Range('D32').Style = $orig_D32
$ws.Range('E32').Value = '  -2.89%  '
$orig_D33 = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.139'
$ws.Range('D33').Style = $orig_D33
$ws.Range('E33').Value = '  -2.51%  '
$orig_D34 = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7033'
$ws.Range('D34').Style = $orig_D34
$ws.Range('E34').Value = '  -4.82%  '
$orig_D35 = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.690'
$ws.Range('D35').Style = $orig_D35
$ws.Range('E35').Value = '  -0.87%  '
$ws.Range('E36').Value = '  -3.61%  '
$orig_D37 = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.647'
$ws.Range('D37').Style = $orig_D37
$ws.Range('E37').Value = '  +0.44%  '
$orig_D38 = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.8867'
$ws.Range('D38').Style = $orig_D38
$ws.Range('E38').Value = '  -2.80%  '
$orig_D39 = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '104.68'
$ws.Range('D39').Style = $orig_D39
$ws.Range('E39').Value = '  -1.58%  '
$orig_D40 = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.919'
$ws.Range('D40').Style = $orig_D40
$ws.Range('E40').Value = '  -5.75%  '
$ws.Range('E41').Value = '  +0.55%  '
$orig_D42 = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.525'
$ws.Range('D42').Style = $orig_D42
$ws.Range('E42').Value = '  -0.67%  '
$orig_D43 = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4009'
$ws.Range('D43').Style = $orig_D43
$ws.Range('E44').Value = '  -1.69%  '
$ws.Range('E45').Value = '  -2.17%  '
$orig_D46 = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '59.93'
$ws.Range('D46').Style = $orig_D46
$ws.Range('E46').Value = '  -2.69%  '
$orig_D47 = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.596'
$ws.Range('D47').Style = $orig_D47
$ws.Range('E47').Value = '  -3.22%  '
$orig_D48 = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '33.07'
$ws.Range('D48').Style = $orig_D48
$ws.Range('E48').Value = '  -1.34%  '
$orig_D49 = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05506'
$ws.Range('D49').Style = $orig_D49
$ws.Range('E49').Value = '  -2.32%  '
$orig_D50 = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.349'
$ws.Range('D50').Style = $orig_D50
$ws.Range('E50').Value = '  -4.36%  '
$orig_D51 = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3630'
$ws.Range('D51').Style = $orig_D51
$ws.Range('E51').Value = '  -3.54%  '
